# Updates cryptos list figures (price + 1h volume change) and re-syncs the
# coin ranking order for rows 21-51 where several coins swapped positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.197.71"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").Value = "1.883.55"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.53%  "

# Row 5
$ws.Range("D5").Value = "'236.98"
$ws.Range("E5").Value = "  -0.73%  "

# Row 6
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.53%  "

# Row 7
$ws.Range("D7").Value = "'0.4675"
$ws.Range("E7").Value = "  -1.94%  "

# Row 8
$ws.Range("D8").Value = "'0.2791"
$ws.Range("E8").Value = "  -2.68%  "

# Row 9
$ws.Range("D9").Value = "'0.06549"
$ws.Range("E9").Value = "  -2.36%  "

# Row 10
$ws.Range("D10").Value = "'19.06"
$ws.Range("E10").Value = "  +1.27%  "

# Row 11
$ws.Range("D11").Value = "'0.07760"
$ws.Range("E11").Value = "  +0.49%  "

# Row 12
$ws.Range("D12").Value = "'97.50"
$ws.Range("E12").Value = "  -5.34%  "

# Row 13
$ws.Range("D13").Value = "1.894.78"
$ws.Range("E13").Value = "  -0.69%  "

# Row 14
$ws.Range("D14").Value = "'5.106"
$ws.Range("E14").Value = "  -2.07%  "

# Row 15
$ws.Range("D15").Value = "'0.6565"
$ws.Range("E15").Value = "  -3.12%  "

# Row 16
$ws.Range("D16").Value = "'279.05"
$ws.Range("E16").Value = "  +6.79%  "

# Row 17
$ws.Range("D17").Value = "30.187.22"
$ws.Range("E17").Value = "  -0.84%  "

# Row 18
$ws.Range("D18").Value = "'1.005"
$ws.Range("E18").Value = "  +0.53%  "

# Row 19
$ws.Range("D19").Value = "2.151.00"
$ws.Range("E19").Value = "  -0.15%  "

# Row 20
$ws.Range("D20").Value = "'12.42"
$ws.Range("E20").Value = "  -2.47%  "

# Row 21
$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").Value = "'1.007"
$ws.Range("E21").Value = "  +0.79%  "

# Row 22
$ws.Range("D22").Value = "'5.310"
$ws.Range("E22").Value = "  -2.03%  "

# Row 23
$ws.Range("B23").Value = "ShibaInu"
$ws.Range("C23").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D23").Value = "'0.000007215"
$ws.Range("E23").Value = "  -3.95%  "

# Row 24
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'6.158"
$ws.Range("E24").Value = "  -2.79%  "

# Row 25
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.200"
$ws.Range("E25").Value = "  -2.97%  "

# Row 26
$ws.Range("D26").Value = "'165.43"
$ws.Range("E26").Value = "  +0.53%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.84"
$ws.Range("E27").Value = "  -1.08%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.990"
$ws.Range("E28").Value = "  -3.39%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.384"
$ws.Range("E29").Value = "  +0.59%  "

# Row 30
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.09806"
$ws.Range("E30").Value = "  -3.28%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.463"
$ws.Range("E31").Value = "  -3.88%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.491"
$ws.Range("E32").Value = "  -1.35%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.176"
$ws.Range("E33").Value = "  -2.09%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04647"
$ws.Range("E34").Value = "  -2.97%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7062"
$ws.Range("E35").Value = "  -3.60%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.087"
$ws.Range("E36").Value = "  -2.48%  "

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.715"
$ws.Range("E37").Value = "  +0.23%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01855"
$ws.Range("E38").Value = "  -3.80%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'6.679"
$ws.Range("E39").Value = "  +6.95%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.514"
$ws.Range("E40").Value = "  -3.25%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'72.04"
$ws.Range("E41").Value = "  -4.09%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8673"
$ws.Range("E42").Value = "  -0.10%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'1.931"
$ws.Range("E43").Value = "  -3.27%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.006"
$ws.Range("E44").Value = "  +0.65%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'103.63"
$ws.Range("E45").Value = "  -3.00%  "

# Row 46
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4158"
$ws.Range("E46").Value = "  -2.55%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'978.12"
$ws.Range("E47").Value = "  -6.81%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.209"
$ws.Range("E48").Value = "  -3.82%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.079"
$ws.Range("E49").Value = "  +1.97%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1162"
$ws.Range("E50").Value = "  -3.31%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05752"
$ws.Range("E51").Value = "  +0.80%  "
